# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" and "全部类型" sheets, for the three affected events.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 3, 5, 6
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 789
$wsExpo.Range("F5").Value = 868
$wsExpo.Range("F6").Value = 2136

# Sheet "全部类型" - rows 3, 7, 8 (same events, different row positions)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 789
$wsAll.Range("F7").Value = 868
$wsAll.Range("F8").Value = 2136
